# Auto-generated Word COM-interop script.
# Applies the '#7187' wording correction to 2_SC_001.docx:
# splits several single-run sentences into multiple runs while
# switching verbs from 1st person to 3rd person (and vice versa
# for the 'Nel caso in cui...' paragraph), and relocates the
# stray _GoBack bookmark into the '5.4' paragraph.
$d = $word.ActiveDocument

# --- Case 1: 'Non posso utilizzare l’app' ---
$rng0 = $d.Content
$rng0.Find.Execute("Non posso utilizzare l’app", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start0 = $rng0.Start
$rng0.Text = "Non può utilizzare l’app"
$p0_0 = $d.Range($start0 + 0, $start0 + 8)
$p0_0.Bold = 1
$p0_0.Bold = 0

# --- Case 2: 'Visualizzo il messaggio “Internet assente”' ---
$rng1 = $d.Content
$rng1.Find.Execute("Visualizzo il messaggio “Internet assente”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start1 = $rng1.Start
$rng1.Text = "Visualizza il messaggio “Internet assente”"
$p1_0 = $d.Range($start1 + 0, $start1 + 10)
$p1_0.Bold = 1
$p1_0.Bold = 0

# --- Case 3: '4.1 Non posso avviare contenuti multimediali' ---
$rng2 = $d.Content
$rng2.Find.Execute("4.1 Non posso avviare contenuti multimediali", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $rng2.Start
$rng2.Text = "4.1 Non può avviare contenuti multimediali"
$p2_0 = $d.Range($start2 + 0, $start2 + 11)
$p2_0.Bold = 1
$p2_0.Bold = 0

# --- Case 4: '4.2 Visualizzo il messaggio “Cuffie non inserite”' ---
$rng3 = $d.Content
$rng3.Find.Execute("4.2 Visualizzo il messaggio “Cuffie non inserite”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start3 = $rng3.Start
$rng3.Text = "4.2 Visualizza il messaggio “Cuffie non inserite”"
$p3_0 = $d.Range($start3 + 0, $start3 + 14)
$p3_0.Bold = 1
$p3_0.Bold = 0

# --- Case 5: '5.4 Non visualizzo contenuti multimediali' ---
$rng4 = $d.Content
$rng4.Find.Execute("5.4 Non visualizzo contenuti multimediali", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start4 = $rng4.Start
$rng4.Text = "5.4 Non visualizza contenuti multimediali"
$p4_0 = $d.Range($start4 + 0, $start4 + 4)
$p4_0.Bold = 1
$p4_0.Bold = 0
$d.Bookmarks.Add("_GoBack", $d.Range($start4 + 4, $start4 + 4)) | Out-Null
$p4_1 = $d.Range($start4 + 4, $start4 + 18)
$p4_1.Bold = 1
$p4_1.Bold = 0

# --- Case 6: '5.5 Visualizzo il messaggio “Contenuti multimediali non disponibili”' ---
$rng5 = $d.Content
$rng5.Find.Execute("5.5 Visualizzo il messaggio “Contenuti multimediali non disponibili”", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start5 = $rng5.Start
$rng5.Text = "5.5 Visualizza il messaggio “Contenuti multimediali non disponibili”"
$p5_0 = $d.Range($start5 + 0, $start5 + 14)
$p5_0.Bold = 1
$p5_0.Bold = 0

# --- Case 7: 'Il dipendente nel caso in cui inserisce username/password errate e si trova al punto 4:' ---
$rng6 = $d.Content
$rng6.Find.Execute("Il dipendente nel caso in cui inserisce username/password errate e si trova al punto 4:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start6 = $rng6.Start
$rng6.Text = "Nel caso in cui il dipendente inserisce username/password errate e siamo al punto 4:"
$p6_0 = $d.Range($start6 + 0, $start6 + 11)
$p6_0.Bold = 1
$p6_0.Bold = 0
$p6_1 = $d.Range($start6 + 11, $start6 + 16)
$p6_1.Bold = 1
$p6_1.Bold = 0
$p6_2 = $d.Range($start6 + 16, $start6 + 30)
$p6_2.Bold = 1
$p6_2.Bold = 0
$p6_3 = $d.Range($start6 + 30, $start6 + 38)
$p6_3.Bold = 1
$p6_3.Bold = 0
$p6_4 = $d.Range($start6 + 38, $start6 + 72)
$p6_4.Bold = 1
$p6_4.Bold = 0

Write-Output "edit.ps1 completed successfully"
